$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 884.875
$ws.Range("J17").Value = 884.875
$ws.Range("L17").Value = 2654.625
$ws.Range("N17").Value = -2990.625

$ws.Range("H108").Value = 45434.668
$ws.Range("J108").Value = 45434.668
$ws.Range("L108").Value = 45434.668
$ws.Range("N108").Value = -53114.668

$ws.Range("H109").Value = 36647.5
$ws.Range("J109").Value = 36647.5
$ws.Range("L109").Value = 36647.5
$ws.Range("N109").Value = -39421.5

$ws.Range("H113").Value = 2187.8262
$ws.Range("I113").Value = 1825
$ws.Range("J113").Value = 2242.25
$ws.Range("K113").Value = 1825
$ws.Range("L113").Value = 2242.25
$ws.Range("M113").Value = 1429
$ws.Range("N113").Value = -8750.25

$ws.Range("H124").Value = 48082.75
$ws.Range("J124").Value = 48082.75
$ws.Range("L124").Value = 48082.75
$ws.Range("N124").Value = -57902.75

$ws.Range("H125").Value = 848.5
$ws.Range("I125").Value = 1461.3334
$ws.Range("J125").Value = 644.2222
$ws.Range("K125").Value = 13152.0006
$ws.Range("L125").Value = 5797.999800000001
$ws.Range("M125").Value = -10692.0006
$ws.Range("N125").Value = -10717.9998

$ws.Range("H128").Value = 40663.2
$ws.Range("J128").Value = 40663.2
$ws.Range("L128").Value = 40663.2
$ws.Range("N128").Value = -50623.2

$ws.Range("H130").Value = 49776
$ws.Range("J130").Value = 49776
$ws.Range("L130").Value = 49776
$ws.Range("N130").Value = -59816

$ws.Range("H138").Value = 2529.967
$ws.Range("I138").Value = 2658.1482
$ws.Range("J138").Value = 2475.8906
$ws.Range("K138").Value = 7974.444600000001
$ws.Range("L138").Value = 7427.6718
$ws.Range("M138").Value = -2834.444600000001
$ws.Range("N138").Value = -17707.6718

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11197.558
$ws.Range("I32").Value = 10303.34
$ws.Range("K32").Value = 10303.34
$ws.Range("M32").Value = -10016.34

$ws.Range("H74").Value = 1438.5818
$ws.Range("I74").Value = 1399.0488
$ws.Range("J74").Value = 1554.3572
$ws.Range("K74").Value = 1399.0488
$ws.Range("L74").Value = 1554.3572
$ws.Range("M74").Value = -525.0488
$ws.Range("N74").Value = -3302.3572

$ws.Range("H77").Value = 1438.5818
$ws.Range("I77").Value = 1399.0488
$ws.Range("J77").Value = 1554.3572
$ws.Range("K77").Value = 6995.244000000001
$ws.Range("L77").Value = 7771.786
$ws.Range("M77").Value = -2627.244000000001
$ws.Range("N77").Value = -16507.786

$ws.Range("H107").Value = 38995
$ws.Range("J107").Value = 38995
$ws.Range("L107").Value = 38995
$ws.Range("N107").Value = -46675

$ws.Range("H109").Value = 40420.25
$ws.Range("J109").Value = 40420.25
$ws.Range("L109").Value = 40420.25
$ws.Range("N109").Value = -43194.25

$ws.Range("H111").Value = 49620
$ws.Range("J111").Value = 49620
$ws.Range("L111").Value = 49620
$ws.Range("N111").Value = -57800

$ws.Range("H117").Value = 40192
$ws.Range("J117").Value = 40192
$ws.Range("L117").Value = 40192
$ws.Range("N117").Value = -49370

$ws.Range("H118").Value = 49409
$ws.Range("J118").Value = 49409
$ws.Range("L118").Value = 49409
$ws.Range("N118").Value = -52723

$ws.Range("H120").Value = 45706
$ws.Range("J120").Value = 45706
$ws.Range("L120").Value = 45706
$ws.Range("N120").Value = -55382

$ws.Range("H125").Value = 44807.332
$ws.Range("J125").Value = 44807.332
$ws.Range("L125").Value = 44807.332
$ws.Range("N125").Value = -54647.332

$ws.Range("H131").Value = 47368.668
$ws.Range("J131").Value = 47368.668
$ws.Range("L131").Value = 47368.668
$ws.Range("N131").Value = -57448.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 47676
$ws.Range("J108").Value = 47676
$ws.Range("L108").Value = 47676
$ws.Range("N108").Value = -55356

$ws.Range("H110").Value = 45231.332
$ws.Range("J110").Value = 45231.332
$ws.Range("L110").Value = 45231.332
$ws.Range("N110").Value = -53411.332

$ws.Range("H111").Value = 47694
$ws.Range("J111").Value = 47694
$ws.Range("L111").Value = 47694
$ws.Range("N111").Value = -55874

$ws.Range("H116").Value = 44630
$ws.Range("J116").Value = 44630
$ws.Range("L116").Value = 44630
$ws.Range("N116").Value = -53808

$ws.Range("H117").Value = 49734
$ws.Range("J117").Value = 49734
$ws.Range("L117").Value = 49734
$ws.Range("N117").Value = -58912

$ws.Range("H130").Value = 47664.4
$ws.Range("J130").Value = 47664.4
$ws.Range("L130").Value = 47664.4
$ws.Range("N130").Value = -57704.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49356.4
$ws.Range("J20").Value = 49356.4
$ws.Range("L20").Value = 49356.4
$ws.Range("N20").Value = -49828.4

$ws.Range("H30").Value = 49356.4
$ws.Range("J30").Value = 49356.4
$ws.Range("L30").Value = 49356.4
$ws.Range("N30").Value = -49538.4

$ws.Range("H58").Value = 1811.8334
$ws.Range("I58").Value = 1559.5555
$ws.Range("J58").Value = 2265.9333
$ws.Range("K58").Value = 1559.5555
$ws.Range("L58").Value = 2265.9333
$ws.Range("M58").Value = -1356.5555
$ws.Range("N58").Value = -2671.9333

$ws.Range("H110").Value = 44701
$ws.Range("J110").Value = 44701
$ws.Range("L110").Value = 44701
$ws.Range("N110").Value = -52881

$ws.Range("H112").Value = 39897.332
$ws.Range("J112").Value = 39897.332
$ws.Range("L112").Value = 39897.332
$ws.Range("N112").Value = -42851.332

$ws.Range("H116").Value = 48489
$ws.Range("J116").Value = 48489
$ws.Range("L116").Value = 48489
$ws.Range("N116").Value = -57667

$ws.Range("H119").Value = 42530
$ws.Range("J119").Value = 42530
$ws.Range("L119").Value = 42530
$ws.Range("N119").Value = -52206

$ws.Range("H128").Value = 49356.4
$ws.Range("J128").Value = 49356.4
$ws.Range("L128").Value = 49356.4
$ws.Range("N128").Value = -59316.4

$ws.Range("H136").Value = 1811.8334
$ws.Range("I136").Value = 1559.5555
$ws.Range("J136").Value = 2265.9333
$ws.Range("K136").Value = 4678.666499999999
$ws.Range("L136").Value = 6797.7999
$ws.Range("M136").Value = -2128.666499999999
$ws.Range("N136").Value = -11897.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 50900
$ws.Range("J110").Value = 50900
$ws.Range("L110").Value = 50900
$ws.Range("N110").Value = -59080

$ws.Range("H114").Value = 38125.668
$ws.Range("J114").Value = 38125.668
$ws.Range("L114").Value = 38125.668
$ws.Range("N114").Value = -46803.668

$ws.Range("H119").Value = 48761
$ws.Range("J119").Value = 48761
$ws.Range("L119").Value = 48761
$ws.Range("N119").Value = -58437

$ws.Range("H130").Value = 44724.668
$ws.Range("J130").Value = 44724.668
$ws.Range("L130").Value = 44724.668
$ws.Range("N130").Value = -54764.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 42939.332
$ws.Range("J118").Value = 42939.332
$ws.Range("L118").Value = 42939.332
$ws.Range("N118").Value = -46253.332

$ws.Range("H124").Value = 46872
$ws.Range("J124").Value = 46872
$ws.Range("L124").Value = 46872
$ws.Range("N124").Value = -56692

$ws.Range("H125").Value = 48711
$ws.Range("J125").Value = 48711
$ws.Range("L125").Value = 48711
$ws.Range("N125").Value = -58551

$ws.Range("H127").Value = 50709
$ws.Range("J127").Value = 50709
$ws.Range("L127").Value = 50709
$ws.Range("N127").Value = -60629

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H130").Value = 35000
$ws.Range("J130").Value = 35000
$ws.Range("L130").Value = 35000
$ws.Range("N130").Value = -45040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 46000
$ws.Range("J16").Value = 46000
$ws.Range("L16").Value = 46000
$ws.Range("N16").Value = -46584

$ws.Range("H108").Value = 46626
$ws.Range("J108").Value = 46626
$ws.Range("L108").Value = 46626
$ws.Range("N108").Value = -54306

$ws.Range("H110").Value = 48644
$ws.Range("J110").Value = 48644
$ws.Range("L110").Value = 48644
$ws.Range("N110").Value = -56824

$ws.Range("H116").Value = 49672
$ws.Range("J116").Value = 49672
$ws.Range("L116").Value = 49672
$ws.Range("N116").Value = -58850

$ws.Range("H128").Value = 48445.5
$ws.Range("J128").Value = 48445.5
$ws.Range("L128").Value = 48445.5
$ws.Range("N128").Value = -58405.5

$ws.Range("H131").Value = 50584
$ws.Range("J131").Value = 50584
$ws.Range("L131").Value = 50584
$ws.Range("N131").Value = -60664
